$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 4.063713
$ws.Range("H2").Value = 12.191139
$ws.Range("I2").Value = 0.5065008440615062
$ws.Range("J2").Value = 0.5065008440615063
$ws.Range("M2").Value = 7.162265000000001
$ws.Range("N2").Value = 21.486795
$ws.Range("O2").Value = 0.4597163377432319
$ws.Range("P2").Value = 0.4597163377432319
$ws.Range("Q2").Value = 29.105389389945
$ws.Range("R2").Value = 261.948504509505
$ws.Range("S2").Value = 0.2328467130958114
$ws.Range("T2").Value = 0.2328467130958115
$ws.Range("G3").Value = 4.063713
$ws.Range("H3").Value = 12.191139
$ws.Range("I3").Value = 0.5065008440615062
$ws.Range("J3").Value = 0.5065008440615063
$ws.Range("O3").Value = 0.2010636628584039
$ws.Range("P3").Value = 0.2010636628584039
$ws.Range("Q3").Value = 12.729667665044
$ws.Range("R3").Value = 114.567008985396
$ws.Range("S3").Value = 0.1018389149478797
$ws.Range("T3").Value = 0.1018389149478797
$ws.Range("G4").Value = 4.063713
$ws.Range("H4").Value = 12.191139
$ws.Range("I4").Value = 0.5065008440615062
$ws.Range("J4").Value = 0.5065008440615063
$ws.Range("M4").Value = 4.518526666666666
$ws.Range("N4").Value = 13.55558
$ws.Range("O4").Value = 0.2900256456854267
$ws.Range("P4").Value = 0.2900256456854267
$ws.Range("Q4").Value = 18.36199555618
$ws.Range("R4").Value = 165.25796000562
$ws.Range("S4").Value = 0.146898234339152
$ws.Range("T4").Value = 0.146898234339152
$ws.Range("G5").Value = 4.063713
$ws.Range("H5").Value = 12.191139
$ws.Range("I5").Value = 0.5065008440615062
$ws.Range("J5").Value = 0.5065008440615063
$ws.Range("M5").Value = 0.7664356666666667
$ws.Range("N5").Value = 2.299307
$ws.Range("O5").Value = 0.04919435371293752
$ws.Range("P5").Value = 0.04919435371293752
$ws.Range("Q5").Value = 3.114574582297
$ws.Range("R5").Value = 28.031171240673
$ws.Range("S5").Value = 0.02491698167866314
$ws.Range("T5").Value = 0.02491698167866315
$ws.Range("I6").Value = 0.2604012840237886
$ws.Range("J6").Value = 0.2604012840237886
$ws.Range("M6").Value = 7.162265000000001
$ws.Range("N6").Value = 21.486795
$ws.Range("O6").Value = 0.4597163377432319
$ws.Range("P6").Value = 0.4597163377432319
$ws.Range("Q6").Value = 14.96360935626333
$ws.Range("R6").Value = 134.67248420637
$ws.Range("S6").Value = 0.1197107246350513
$ws.Range("T6").Value = 0.1197107246350513
$ws.Range("I7").Value = 0.2604012840237886
$ws.Range("J7").Value = 0.2604012840237886
$ws.Range("O7").Value = 0.2010636628584039
$ws.Range("P7").Value = 0.2010636628584039
$ws.Range("S7").Value = 0.05235723597885449
$ws.Range("T7").Value = 0.05235723597885451
$ws.Range("I8").Value = 0.2604012840237886
$ws.Range("J8").Value = 0.2604012840237886
$ws.Range("M8").Value = 4.518526666666666
$ws.Range("N8").Value = 13.55558
$ws.Range("O8").Value = 0.2900256456854267
$ws.Range("P8").Value = 0.2900256456854267
$ws.Range("Q8").Value = 9.440235443097777
$ws.Range("R8").Value = 84.96211898788
$ws.Range("S8").Value = 0.07552305053631347
$ws.Range("T8").Value = 0.07552305053631349
$ws.Range("I9").Value = 0.2604012840237886
$ws.Range("J9").Value = 0.2604012840237886
$ws.Range("M9").Value = 0.7664356666666667
$ws.Range("N9").Value = 2.299307
$ws.Range("O9").Value = 0.04919435371293752
$ws.Range("P9").Value = 0.04919435371293752
$ws.Range("Q9").Value = 1.601259365955778
$ws.Range("R9").Value = 14.411334293602
$ws.Range("S9").Value = 0.01281027287356936
$ws.Range("T9").Value = 0.01281027287356936
$ws.Range("G10").Value = 1.588356333333333
$ws.Range("H10").Value = 4.765069
$ws.Range("I10").Value = 0.1979725988286506
$ws.Range("J10").Value = 0.1979725988286507
$ws.Range("M10").Value = 7.162265000000001
$ws.Range("N10").Value = 21.486795
$ws.Range("O10").Value = 0.4597163377432319
$ws.Range("P10").Value = 0.4597163377432319
$ws.Range("Q10").Value = 11.37622897376167
$ws.Range("R10").Value = 102.386060763855
$ws.Range("S10").Value = 0.09101123810701732
$ws.Range("T10").Value = 0.09101123810701733
$ws.Range("G11").Value = 1.588356333333333
$ws.Range("H11").Value = 4.765069
$ws.Range("I11").Value = 0.1979725988286506
$ws.Range("J11").Value = 0.1979725988286507
$ws.Range("O11").Value = 0.2010636628584039
$ws.Range("P11").Value = 0.2010636628584039
$ws.Range("Q11").Value = 4.975560099101778
$ws.Range("R11").Value = 44.780040891916
$ws.Range("S11").Value = 0.03980509586608585
$ws.Range("T11").Value = 0.03980509586608585
$ws.Range("G12").Value = 1.588356333333333
$ws.Range("H12").Value = 4.765069
$ws.Range("I12").Value = 0.1979725988286506
$ws.Range("J12").Value = 0.1979725988286507
$ws.Range("M12").Value = 4.518526666666666
$ws.Range("N12").Value = 13.55558
$ws.Range("O12").Value = 0.2900256456854267
$ws.Range("P12").Value = 0.2900256456854267
$ws.Range("Q12").Value = 7.177030448335555
$ws.Range("R12").Value = 64.59327403501999
$ws.Range("S12").Value = 0.05741713080330136
$ws.Range("T12").Value = 0.05741713080330137
$ws.Range("G13").Value = 1.588356333333333
$ws.Range("H13").Value = 4.765069
$ws.Range("I13").Value = 0.1979725988286506
$ws.Range("J13").Value = 0.1979725988286507
$ws.Range("M13").Value = 0.7664356666666667
$ws.Range("N13").Value = 2.299307
$ws.Range("O13").Value = 0.04919435371293752
$ws.Range("P13").Value = 0.04919435371293752
$ws.Range("Q13").Value = 1.217372945242556
$ws.Range("R13").Value = 10.956356507183
$ws.Range("S13").Value = 0.009739134052246119
$ws.Range("T13").Value = 0.009739134052246121
$ws.Range("G14").Value = 0.281814
$ws.Range("H14").Value = 0.845442
$ws.Range("I14").Value = 0.03512527308605438
$ws.Range("J14").Value = 0.03512527308605439
$ws.Range("M14").Value = 7.162265000000001
$ws.Range("N14").Value = 21.486795
$ws.Range("O14").Value = 0.4597163377432319
$ws.Range("P14").Value = 0.4597163377432319
$ws.Range("Q14").Value = 2.01842654871
$ws.Range("R14").Value = 18.16583893839
$ws.Range("S14").Value = 0.01614766190535183
$ws.Range("T14").Value = 0.01614766190535183
$ws.Range("G15").Value = 0.281814
$ws.Range("H15").Value = 0.845442
$ws.Range("I15").Value = 0.03512527308605438
$ws.Range("J15").Value = 0.03512527308605439
$ws.Range("O15").Value = 0.2010636628584039
$ws.Range("P15").Value = 0.2010636628584039
$ws.Range("Q15").Value = 0.882788367032
$ws.Range("R15").Value = 7.945095303288001
$ws.Range("S15").Value = 0.007062416065583805
$ws.Range("T15").Value = 0.007062416065583807
$ws.Range("G16").Value = 0.281814
$ws.Range("H16").Value = 0.845442
$ws.Range("I16").Value = 0.03512527308605438
$ws.Range("J16").Value = 0.03512527308605439
$ws.Range("M16").Value = 4.518526666666666
$ws.Range("N16").Value = 13.55558
$ws.Range("O16").Value = 0.2900256456854267
$ws.Range("P16").Value = 0.2900256456854267
$ws.Range("Q16").Value = 1.27338407404
$ws.Range("R16").Value = 11.46045666636
$ws.Range("S16").Value = 0.01018723000665986
$ws.Range("T16").Value = 0.01018723000665987
$ws.Range("G17").Value = 0.281814
$ws.Range("H17").Value = 0.845442
$ws.Range("I17").Value = 0.03512527308605438
$ws.Range("J17").Value = 0.03512527308605439
$ws.Range("M17").Value = 0.7664356666666667
$ws.Range("N17").Value = 2.299307
$ws.Range("O17").Value = 0.04919435371293752
$ws.Range("P17").Value = 0.04919435371293752
$ws.Range("Q17").Value = 0.215992300966
$ws.Range("R17").Value = 1.943930708694
$ws.Range("S17").Value = 0.001727965108458884
$ws.Range("T17").Value = 0.001727965108458884
